$d = $word.ActiveDocument

# Locate the paragraph that currently reads "Link:" (it is the final
# paragraph in the document body and its run carries the
# <w:lastRenderedPageBreak/> marker).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Text.Trim() -eq "Link:") {
    $target = $p
  }
}

if ($target -ne $null) {
  $pr = $target.Range
  # Range covering just the paragraph's own content (exclude the
  # end-of-paragraph mark) so the replacement below does not swallow the
  # following section break / extra paragraphs.
  $contentRange = $d.Range($pr.Start, $pr.End - 1)

  # Build the five replacement paragraphs:
  #   1) "22363276 BZS Khumalo"  (keeps the lastRenderedPageBreak marker
  #      that used to sit on the "Link:" run)
  #   2) "22328828 N Cele"
  #   3) "22329111 W Khuzwayo "
  #   4) an empty paragraph
  #   5) "Link: https://fb.watch/sRKk_aNzti/"
  $body = ""
  $body += "<w:p><w:r><w:lastRenderedPageBreak/><w:t>22363276</w:t></w:r>"
  $body += "<w:r><w:t xml:space=`"preserve`"> BZS Khumalo</w:t></w:r></w:p>"
  $body += "<w:p><w:r><w:t>22328828 N Cele</w:t></w:r></w:p>"
  $body += "<w:p><w:r><w:t>22329111</w:t></w:r>"
  $body += "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"
  $body += "<w:r><w:t xml:space=`"preserve`">W Khuzwayo </w:t></w:r></w:p>"
  $body += "<w:p/>"
  $body += "<w:p><w:r><w:t>Link:</w:t></w:r>"
  $body += "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"
  $body += "<w:r><w:t>https://fb.watch/sRKk_aNzti/</w:t></w:r></w:p>"

  $xml = "<?xml version=`"1.0`" standalone=`"yes`"?>"
  $xml += "<?mso-application progid=`"Word.Document`"?>"
  $xml += "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`">"
  $xml += "<pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`">"
  $xml += "<pkg:xmlData>"
  $xml += "<w:document xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`">"
  $xml += "<w:body>" + $body + "</w:body>"
  $xml += "</w:document>"
  $xml += "</pkg:xmlData></pkg:part></pkg:package>"

  $contentRange.InsertXML($xml)
}
